$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.216187119483948
$ws.Range("B1").Value = 2.450452327728271
$ws.Range("C1").Value = 7.247153759002686
$ws.Range("D1").Value = 2.255217552185059
$ws.Range("E1").Value = 1.163832068443298
